$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 401.69232
$ws.Range("I41").Value = 346.25
$ws.Range("J41").Value = 490.4
$ws.Range("K41").Value = 346.25
$ws.Range("L41").Value = 490.4
$ws.Range("M41").Value = 93.75
$ws.Range("N41").Value = -1370.4
$ws.Range("H76").Value = 4395.4614
$ws.Range("I76").Value = 3477.2856
$ws.Range("K76").Value = 3477.2856
$ws.Range("M76").Value = -3162.2856
$ws.Range("H79").Value = 4395.4614
$ws.Range("I79").Value = 3477.2856
$ws.Range("K79").Value = 3477.2856
$ws.Range("M79").Value = -2385.2856
$ws.Range("H111").Value = 817.72974
$ws.Range("J111").Value = 942.6429000000001
$ws.Range("L111").Value = 2827.9287
$ws.Range("N111").Value = -8961.9287
$ws.Range("H113").Value = 8076.25
$ws.Range("I113").Value = 21515
$ws.Range("J113").Value = 3596.6667
$ws.Range("K113").Value = 21515
$ws.Range("L113").Value = 3596.6667
$ws.Range("M113").Value = -18261
$ws.Range("N113").Value = -10104.6667
$ws.Range("H116").Value = 5501.4165
$ws.Range("I116").Value = 3226.25
$ws.Range("J116").Value = 6639
$ws.Range("K116").Value = 3226.25
$ws.Range("L116").Value = 6639
$ws.Range("M116").Value = 215.75
$ws.Range("N116").Value = -13523
$ws.Range("H135").Value = 5881.304
$ws.Range("I135").Value = 7333.875
$ws.Range("J135").Value = 2561.1428
$ws.Range("K135").Value = 66004.875
$ws.Range("L135").Value = 23050.2852
$ws.Range("M135").Value = -63469.875
$ws.Range("N135").Value = -28120.2852
$ws.Range("H137").Value = 239810.56
$ws.Range("I137").Value = 314660
$ws.Range("J137").Value = 36647.785
$ws.Range("K137").Value = 943980
$ws.Range("L137").Value = 109943.355
$ws.Range("M137").Value = -941430
$ws.Range("N137").Value = -115043.355

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1102.75
$ws.Range("I2").Value = 1137
$ws.Range("K2").Value = 1137
$ws.Range("M2").Value = -1024
$ws.Range("H32").Value = 773654.25
$ws.Range("I32").Value = 4109.5454
$ws.Range("K32").Value = 4109.5454
$ws.Range("M32").Value = -3822.5454
$ws.Range("H116").Value = 1102.75
$ws.Range("I116").Value = 1137
$ws.Range("K116").Value = 1137
$ws.Range("M116").Value = 1157
$ws.Range("H125").Value = 49986.668
$ws.Range("J125").Value = 49986.668
$ws.Range("L125").Value = 49986.668
$ws.Range("N125").Value = -59826.668
$ws.Range("H132").Value = 3379621.5
$ws.Range("I132").Value = 3677382.2
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 11032146.6
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -11029616.6
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1102.75
$ws.Range("I3").Value = 1137
$ws.Range("K3").Value = 1137
$ws.Range("M3").Value = -1023
$ws.Range("H134").Value = 8334879
$ws.Range("I134").Value = 11906137
$ws.Range("J134").Value = 1943.5834
$ws.Range("K134").Value = 35718411
$ws.Range("L134").Value = 5830.7502
$ws.Range("M134").Value = -35715876
$ws.Range("N134").Value = -10900.7502

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 39185.08
$ws.Range("I16").Value = 67360.87
$ws.Range("J16").Value = 763.5454999999999
$ws.Range("K16").Value = 67360.87
$ws.Range("L16").Value = 763.5454999999999
$ws.Range("M16").Value = -67073.87
$ws.Range("N16").Value = -1337.5455
$ws.Range("H31").Value = 6086.525
$ws.Range("I31").Value = 901.75
$ws.Range("J31").Value = 18184.334
$ws.Range("K31").Value = 901.75
$ws.Range("L31").Value = 18184.334
$ws.Range("M31").Value = -606.75
$ws.Range("N31").Value = -18774.334
$ws.Range("H34").Value = 6086.525
$ws.Range("I34").Value = 901.75
$ws.Range("J34").Value = 18184.334
$ws.Range("K34").Value = 901.75
$ws.Range("L34").Value = 18184.334
$ws.Range("M34").Value = -699.75
$ws.Range("N34").Value = -18588.334
$ws.Range("H58").Value = 4934073.5
$ws.Range("I58").Value = 6803656
$ws.Range("K58").Value = 6803656
$ws.Range("M58").Value = -6803453
$ws.Range("H113").Value = 39185.08
$ws.Range("I113").Value = 67360.87
$ws.Range("J113").Value = 763.5454999999999
$ws.Range("K113").Value = 67360.87
$ws.Range("L113").Value = 763.5454999999999
$ws.Range("M113").Value = -65190.87
$ws.Range("N113").Value = -5103.5455
$ws.Range("H136").Value = 4934073.5
$ws.Range("I136").Value = 6803656
$ws.Range("K136").Value = 20410968
$ws.Range("M136").Value = -20408418

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 892.7193
$ws.Range("I132").Value = 600.89746
$ws.Range("J132").Value = 1525
$ws.Range("K132").Value = 5408.07714
$ws.Range("L132").Value = 13725
$ws.Range("M132").Value = -2878.07714
$ws.Range("N132").Value = -18785

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1578.2759
$ws.Range("I113").Value = 1322.8889
$ws.Range("J113").Value = 1996.1818
$ws.Range("K113").Value = 1322.8889
$ws.Range("L113").Value = 1996.1818
$ws.Range("M113").Value = 847.1111000000001
$ws.Range("N113").Value = -6336.1818

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1632.9584
$ws.Range("I61").Value = 1142.2858
$ws.Range("J61").Value = 1835
$ws.Range("K61").Value = 1142.2858
$ws.Range("L61").Value = 1835
$ws.Range("M61").Value = -940.2858000000001
$ws.Range("N61").Value = -2239
$ws.Range("H113").Value = 1632.9584
$ws.Range("I113").Value = 1142.2858
$ws.Range("J113").Value = 1835
$ws.Range("K113").Value = 1142.2858
$ws.Range("L113").Value = 1835
$ws.Range("M113").Value = 1027.7142
$ws.Range("N113").Value = -6175

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 453.8889
$ws.Range("I100").Value = 380.4
$ws.Range("J100").Value = 545.75
$ws.Range("K100").Value = 760.8
$ws.Range("L100").Value = 1091.5
$ws.Range("M100").Value = -219.8
$ws.Range("N100").Value = -2173.5
$ws.Range("H132").Value = 59107170
$ws.Range("I132").Value = 60002020
$ws.Range("J132").Value = 57118610
$ws.Range("K132").Value = 180006060
$ws.Range("L132").Value = 171355830
$ws.Range("M132").Value = -180003530
$ws.Range("N132").Value = -171360890
